$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "SA-HW40.xpc" to "SA"
$ws.Name = "SA"

# Add the new row 16 of averaged intensity data (Gaussian Quadrature scheme export)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.001045738816186
$ws.Range("D16").Value = 0.9763106660148633
$ws.Range("E16").Value = 0.9905882352941177
$ws.Range("F16").Value = 0.9888527775402498
$ws.Range("G16").Value = 1.001045738816186
$ws.Range("H16").Value = 0.9763106660148633
$ws.Range("I16").Value = 0.99
$ws.Range("J16").Value = 0.981764705882353
$ws.Range("K16").Value = 0.9925490388345444
$ws.Range("L16").Value = 0.9819315554119415
$ws.Range("M16").Value = 1.001045738816186
$ws.Range("N16").Value = 0.9834494506544904
$ws.Range("O16").Value = 0.9891993544163542
$ws.Range("P16").Value = 0.9878803397242819

# Match the bold/bordered/centered style used on A3:A15 for the new row's A cell
# (copy the existing format from A15 rather than re-specifying it, to avoid
# introducing a redundant cell style definition)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
